$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C (Socket), shifting C:G to D:H
$ws.Columns.Item(3).Insert()

# Set new header for column C
$ws.Range('C1').Value = 'Price'

# Fill Price column (C) for all data rows with the price list string
$price = '[''199'', ''199'', ''598'', ''206'', ''399'', ''179'', ''279'', ''299'', ''779'', ''308'', ''313'', ''164'', ''371'', ''272'', ''312'', ''139'', ''193'', ''399'', ''151'', ''549'', ''169'', ''259'', ''249'', ''149'', ''134'', ''549'', ''317'', ''78'', ''107'', ''328'', ''363'', ''203'', ''209'', ''111'', ''389'', ''508'', ''298'']'
$ws.Range('C2:C38').Value = $price

# Rows 29-38 data got rewritten/rotated in the source; update Name (B) and
# Socket/Cores/Threads/Operating Frequency/Max Operating Frequency (D:H)

# Row 29: Core i3-10100F
$ws.Range('B29').Value = 'Core i3-10100F'
$ws.Range('D29').Value = 'LGA 1200'
$ws.Range('E29').Value = 'Quad-Core'
$ws.Range('F29').Value = '8'
$ws.Range('G29').Value = '3.6 GHz'
$ws.Range('H29').Value = '4.30 GHz'

# Row 30: Core i3-12100F
$ws.Range('B30').Value = 'Core i3-12100F'
$ws.Range('D30').Value = 'LGA 1700'
$ws.Range('E30').Value = 'Quad-Core'
$ws.Range('F30').Value = '8'
$ws.Range('G30').Value = '3.3 GHz'
$ws.Range('H30').Value = '4.3 GHz'

# Row 31: Core i7-10700K
$ws.Range('B31').Value = 'Core i7-10700K'
$ws.Range('D31').Value = 'LGA 1200'
$ws.Range('E31').Value = '8-Core'
$ws.Range('F31').Value = '16'
$ws.Range('G31').Value = '3.8 GHz'
$ws.Range('H31').Value = '5.10 GHz'

# Row 32: Core i7-12700KF
$ws.Range('B32').Value = 'Core i7-12700KF'
$ws.Range('D32').Value = 'LGA 1700'
$ws.Range('E32').Value = '12-Core (8P+4E)'
$ws.Range('F32').Value = '20'
$ws.Range('G32').Value = 'P-core Base Frequency: 3.6 GHzE-core Base Frequency: 2.7 GHz'
$ws.Range('H32').Value = 'Intel Turbo Boost Max Technology 3.0 Frequency: Up to 5.0 GHzSingle P-core Turbo Frequency: Up to 4.9 GHzSingle E-core Turbo Frequency: Up to 3.8 GHz'

# Row 33: Core i5-10600K
$ws.Range('B33').Value = 'Core i5-10600K'
$ws.Range('D33').Value = 'LGA 1200'
$ws.Range('E33').Value = '6-Core'
$ws.Range('F33').Value = '12'
$ws.Range('G33').Value = '4.1 GHz'
$ws.Range('H33').Value = '4.80 GHz'

# Row 34: Core i5-11600KF
$ws.Range('B34').Value = 'Core i5-11600KF'
$ws.Range('D34').Value = 'LGA 1200'
$ws.Range('E34').Value = '6-Core'
$ws.Range('F34').Value = '12'
$ws.Range('G34').Value = '3.9 GHz'
$ws.Range('H34').Value = '4.9 GHz'

# Row 35: Core i3-10105
$ws.Range('B35').Value = 'Core i3-10105'
$ws.Range('D35').Value = 'LGA 1200'
$ws.Range('E35').Value = 'Quad-Core'
$ws.Range('F35').Value = '8'
$ws.Range('G35').Value = '3.7 GHz'
$ws.Range('H35').Value = '4.4 GHz'

# Row 36: Core i9-10900K
$ws.Range('B36').Value = 'Core i9-10900K'
$ws.Range('D36').Value = 'LGA 1200'
$ws.Range('E36').Value = '10-Core'
$ws.Range('F36').Value = '20'
$ws.Range('G36').Value = '3.7 GHz'
$ws.Range('H36').Value = '5.30 GHz'

# Row 37: Core i9-12900
$ws.Range('B37').Value = 'Core i9-12900'
$ws.Range('D37').Value = 'LGA 1700'
$ws.Range('E37').Value = '16-Core (8P+8E)'
$ws.Range('F37').Value = '24'
$ws.Range('G37').Value = 'P-core Base Frequency: 2.4 GHzE-core Base Frequency: 1.8 GHz'
$ws.Range('H37').Value = 'Intel Turbo Boost Max Technology 3.0 Frequency: Up to 5.1 GHzP-core Turbo Frequency: Up to 5.0 GHzE-core Turbo Frequency: Up to 3.8 GHz'

# Row 38: Core i7-11700KF
$ws.Range('B38').Value = 'Core i7-11700KF'
$ws.Range('D38').Value = 'LGA 1200'
$ws.Range('E38').Value = '8-Core'
$ws.Range('F38').Value = '16'
$ws.Range('G38').Value = '3.6 GHz'
$ws.Range('H38').Value = '5.0 GHz'
